# Auto-generated Excel COM-interop script
# Applies numeric cell updates to the Leve-profit sheets (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR)
# as produced by the scheduled market-data refresh runner.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 17
$ws.Range("H17").Value = 2778.6726
$ws.Range("J17").Value = 2854.9812
$ws.Range("L17").Value = 8564.943600000001
$ws.Range("N17").Value = -8900.943600000001
# Row 42
$ws.Range("H42").Value = 5250
$ws.Range("I42").Value = 500
$ws.Range("J42").Value = 7625
$ws.Range("K42").Value = 1500
$ws.Range("L42").Value = 22875
$ws.Range("M42").Value = -1270
$ws.Range("N42").Value = -23335
# Row 58
$ws.Range("H58").Value = 7339.0586
$ws.Range("J58").Value = 12349.7
$ws.Range("L58").Value = 37049.10000000001
$ws.Range("N58").Value = -37349.10000000001
# Row 80
$ws.Range("H80").Value = 609.8
$ws.Range("J80").Value = 802.5
$ws.Range("L80").Value = 2407.5
$ws.Range("N80").Value = -4403.5
# Row 83
$ws.Range("H83").Value = 609.8
$ws.Range("J83").Value = 802.5
$ws.Range("L83").Value = 7222.5
$ws.Range("N83").Value = -17206.5
# Row 98
$ws.Range("H98").Value = 2636.7144
$ws.Range("I98").Value = 2460.75
$ws.Range("J98").Value = 3199.8
$ws.Range("K98").Value = 2460.75
$ws.Range("L98").Value = 3199.8
$ws.Range("M98").Value = -962.75
$ws.Range("N98").Value = -6195.8
# Row 122
$ws.Range("H122").Value = 2636.7144
$ws.Range("I122").Value = 2460.75
$ws.Range("J122").Value = 3199.8
$ws.Range("K122").Value = 7382.25
$ws.Range("L122").Value = 9599.400000000001
$ws.Range("M122").Value = -4932.25
$ws.Range("N122").Value = -14499.4
# Row 132
$ws.Range("H132").Value = 25641852
$ws.Range("I132").Value = 31250534
$ws.Range("J132").Value = 2154.7144
$ws.Range("K132").Value = 93751602
$ws.Range("L132").Value = 6464.1432
$ws.Range("M132").Value = -93749072
$ws.Range("N132").Value = -11524.1432
# Row 135
$ws.Range("H135").Value = 2418
$ws.Range("I135").Value = 2720.25
$ws.Range("J135").Value = 1209
$ws.Range("K135").Value = 24482.25
$ws.Range("L135").Value = 10881
$ws.Range("M135").Value = -21947.25
$ws.Range("N135").Value = -15951

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 3698.1177
$ws.Range("I32").Value = 2252.652
$ws.Range("K32").Value = 2252.652
$ws.Range("M32").Value = -1965.652
# Row 45
$ws.Range("H45").Value = 1536.7
$ws.Range("I45").Value = 1409.3077
$ws.Range("J45").Value = 1773.2858
$ws.Range("K45").Value = 1409.3077
$ws.Range("L45").Value = 1773.2858
$ws.Range("M45").Value = -1032.3077
$ws.Range("N45").Value = -2527.2858
# Row 61
$ws.Range("H61").Value = 1983.4615
$ws.Range("I61").Value = 1961.25
$ws.Range("K61").Value = 1961.25
$ws.Range("M61").Value = -1749.25
# Row 74
$ws.Range("H74").Value = 2894
$ws.Range("I74").Value = 1852.7858
$ws.Range("J74").Value = 4351.7
$ws.Range("K74").Value = 1852.7858
$ws.Range("L74").Value = 4351.7
$ws.Range("M74").Value = -978.7858000000001
$ws.Range("N74").Value = -6099.7
# Row 77
$ws.Range("H77").Value = 2894
$ws.Range("I77").Value = 1852.7858
$ws.Range("J77").Value = 4351.7
$ws.Range("K77").Value = 9263.929
$ws.Range("L77").Value = 21758.5
$ws.Range("M77").Value = -4895.929
$ws.Range("N77").Value = -30494.5
# Row 101
$ws.Range("H101").Value = 206767
$ws.Range("J101").Value = 206767
$ws.Range("L101").Value = 206767
$ws.Range("N101").Value = -213257
# Row 122
$ws.Range("H122").Value = 6901865
$ws.Range("I122").Value = 7148128
$ws.Range("K122").Value = 21444384
$ws.Range("M122").Value = -21441934
# Row 132
$ws.Range("H132").Value = 5841.352
$ws.Range("I132").Value = 3507.0635
$ws.Range("K132").Value = 10521.1905
$ws.Range("M132").Value = -7991.190500000001
# Row 136
$ws.Range("H136").Value = 1983.4615
$ws.Range("I136").Value = 1961.25
$ws.Range("K136").Value = 5883.75
$ws.Range("M136").Value = -3333.75

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 55556980
$ws.Range("I20").Value = 71429770
$ws.Range("K20").Value = 71429770
$ws.Range("M20").Value = -71429523
# Row 107
$ws.Range("H107").Value = 3922.5293
$ws.Range("I107").Value = 3076.2222
$ws.Range("J107").Value = 4874.625
$ws.Range("K107").Value = 3076.2222
$ws.Range("L107").Value = 4874.625
$ws.Range("M107").Value = -1156.2222
$ws.Range("N107").Value = -8714.625
# Row 112
$ws.Range("H112").Value = 62499.5
$ws.Range("J112").Value = 49999
$ws.Range("L112").Value = 49999
$ws.Range("N112").Value = -52953
# Row 134
$ws.Range("H134").Value = 16668150
$ws.Range("I134").Value = 16668150
$ws.Range("J134").Value = 0
$ws.Range("K134").Value = 50004450
$ws.Range("L134").Value = 0
$ws.Range("M134").Value = -50001915
$ws.Range("N134").ClearContents()
# Row 135
$ws.Range("H135").Value = 65689.67999999999
$ws.Range("J135").Value = 65689.67999999999
$ws.Range("L135").Value = 65689.67999999999
$ws.Range("N135").Value = -75829.67999999999

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 58
$ws.Range("H58").Value = 4282.4287
$ws.Range("I58").Value = 4669.706
$ws.Range("K58").Value = 4669.706
$ws.Range("M58").Value = -4466.706
# Row 122
$ws.Range("H122").Value = 3162.6667
$ws.Range("I122").Value = 2647.25
$ws.Range("J122").Value = 3751.7144
$ws.Range("K122").Value = 7941.75
$ws.Range("L122").Value = 11255.1432
$ws.Range("M122").Value = -5491.75
$ws.Range("N122").Value = -16155.1432
# Row 132
$ws.Range("H132").Value = 3404.5715
$ws.Range("I132").Value = 2969.4546
$ws.Range("J132").Value = 5000
$ws.Range("K132").Value = 8908.363799999999
$ws.Range("L132").Value = 15000
$ws.Range("M132").Value = -6378.363799999999
$ws.Range("N132").Value = -20060
# Row 136
$ws.Range("H136").Value = 4282.4287
$ws.Range("I136").Value = 4669.706
$ws.Range("K136").Value = 14009.118
$ws.Range("M136").Value = -11459.118

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 113
$ws.Range("H113").Value = 13890863
$ws.Range("I113").Value = 2031.1428
$ws.Range("J113").Value = 19609794
$ws.Range("K113").Value = 6093.428400000001
$ws.Range("L113").Value = 58829382
$ws.Range("M113").Value = -3923.428400000001
$ws.Range("N113").Value = -58833722

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 70
$ws.Range("H70").Value = 34025.715
$ws.Range("I70").Value = 50644.082
$ws.Range("J70").Value = 11867.889
$ws.Range("K70").Value = 50644.082
$ws.Range("L70").Value = 11867.889
$ws.Range("M70").Value = -50374.082
$ws.Range("N70").Value = -12407.889
# Row 73
$ws.Range("H73").Value = 34025.715
$ws.Range("I73").Value = 50644.082
$ws.Range("J73").Value = 11867.889
$ws.Range("K73").Value = 50644.082
$ws.Range("L73").Value = 11867.889
$ws.Range("M73").Value = -49708.082
$ws.Range("N73").Value = -13739.889
# Row 80
$ws.Range("H80").Value = 4000.5
$ws.Range("I80").Value = 4102
$ws.Range("J80").Value = 3949.75
$ws.Range("K80").Value = 4102
$ws.Range("L80").Value = 3949.75
$ws.Range("M80").Value = -3104
$ws.Range("N80").Value = -5945.75
# Row 83
$ws.Range("H83").Value = 4000.5
$ws.Range("I83").Value = 4102
$ws.Range("J83").Value = 3949.75
$ws.Range("K83").Value = 20510
$ws.Range("L83").Value = 19748.75
$ws.Range("M83").Value = -15518
$ws.Range("N83").Value = -29732.75
# Row 121
$ws.Range("H121").Value = 40000
$ws.Range("J121").Value = 40000
$ws.Range("L121").Value = 40000
$ws.Range("N121").Value = -43494
# Row 132
$ws.Range("H132").Value = 2090.75
$ws.Range("I132").Value = 1343.3334
$ws.Range("J132").Value = 4333
$ws.Range("K132").Value = 4030.0002
$ws.Range("L132").Value = 12999
$ws.Range("M132").Value = -1500.0002
$ws.Range("N132").Value = -18059

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 16
$ws.Range("H16").Value = 852.1786
$ws.Range("I16").Value = 969.4545000000001
$ws.Range("J16").Value = 422.16666
$ws.Range("K16").Value = 969.4545000000001
$ws.Range("L16").Value = 422.16666
$ws.Range("M16").Value = -799.4545000000001
$ws.Range("N16").Value = -762.16666
# Row 40
$ws.Range("H40").Value = 5156.857
$ws.Range("I40").Value = 4818.1924
$ws.Range("K40").Value = 4818.1924
$ws.Range("M40").Value = -4682.1924
# Row 55
$ws.Range("H55").Value = 994.4545000000001
$ws.Range("I55").Value = 380.53845
$ws.Range("J55").Value = 1881.2222
$ws.Range("K55").Value = 380.53845
$ws.Range("L55").Value = 1881.2222
$ws.Range("M55").Value = -207.53845
$ws.Range("N55").Value = -2227.2222
# Row 68
$ws.Range("H68").Value = 2571.1428
$ws.Range("I68").Value = 2666.3333
$ws.Range("J68").Value = 2000
$ws.Range("K68").Value = 2666.3333
$ws.Range("L68").Value = 2000
$ws.Range("M68").Value = -1917.3333
$ws.Range("N68").Value = -3498
# Row 71
$ws.Range("H71").Value = 2571.1428
$ws.Range("I71").Value = 2666.3333
$ws.Range("J71").Value = 2000
$ws.Range("K71").Value = 13331.6665
$ws.Range("L71").Value = 10000
$ws.Range("M71").Value = -9587.666499999999
$ws.Range("N71").Value = -17488
# Row 82
$ws.Range("H82").Value = 52633330
$ws.Range("J82").Value = 1618.8
$ws.Range("L82").Value = 1618.8
$ws.Range("N82").Value = -2340.8
# Row 85
$ws.Range("H85").Value = 52633330
$ws.Range("J85").Value = 1618.8
$ws.Range("L85").Value = 1618.8
$ws.Range("N85").Value = -4114.8
# Row 122
$ws.Range("H122").Value = 9305.883
$ws.Range("I122").Value = 5561.615
$ws.Range("K122").Value = 16684.845
$ws.Range("M122").Value = -14234.845
# Row 132
$ws.Range("H132").Value = 5750.077
$ws.Range("I132").Value = 4896.0835
$ws.Range("K132").Value = 14688.2505
$ws.Range("M132").Value = -12158.2505
# Row 136
$ws.Range("H136").Value = 3607
$ws.Range("I136").Value = 2672
$ws.Range("K136").Value = 8016
$ws.Range("M136").Value = -5466
# Row 138
$ws.Range("H138").Value = 90000
$ws.Range("J138").Value = 90000
$ws.Range("L138").Value = 90000
$ws.Range("N138").Value = -100280

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 96
$ws.Range("H96").Value = 18518818
$ws.Range("I96").Value = 18518818
$ws.Range("J96").Value = 0
$ws.Range("K96").Value = 18518818
$ws.Range("L96").Value = 0
$ws.Range("M96").Value = -18517445
$ws.Range("N96").ClearContents()
# Row 126
$ws.Range("H126").Value = 32259082
$ws.Range("I126").Value = 969.16
$ws.Range("J126").Value = 166667890
$ws.Range("K126").Value = 2907.48
$ws.Range("L126").Value = 500003670
$ws.Range("M126").Value = -437.48
$ws.Range("N126").Value = -500008610
# Row 132
$ws.Range("H132").Value = 1398.4166
$ws.Range("I132").Value = 1274.2
$ws.Range("J132").Value = 2019.5
$ws.Range("K132").Value = 3822.6
$ws.Range("L132").Value = 6058.5
$ws.Range("M132").Value = -1292.6
$ws.Range("N132").Value = -11118.5

